# Reorder the "Recorded By" names in column G so that "System" (or "system")
# moves from the front of the comma-separated list to the end, keeping the
# other entries in their existing relative order.
#
# Examples of the transformation performed:
#   "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, system, System"     -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    $text = [string]$value

    if ($text -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($text -eq "backup@backdoor.com, system, System") {
        $cell.Value2 = "system, backup@backdoor.com, System"
    }
}
